# Update the "想去人数" (want-to-go count) figures in column F across the
# "展览", "演出" and "全部类型" sheets to match the newly generated data.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" ---
$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F2").Value  = 289
$wsExhibit.Range("F4").Value  = 10197
$wsExhibit.Range("F5").Value  = 330
$wsExhibit.Range("F6").Value  = 932
$wsExhibit.Range("F7").Value  = 1266
$wsExhibit.Range("F8").Value  = 6611
$wsExhibit.Range("F10").Value = 427
$wsExhibit.Range("F12").Value = 124
$wsExhibit.Range("F13").Value = 3145
$wsExhibit.Range("F16").Value = 619
$wsExhibit.Range("F17").Value = 119
$wsExhibit.Range("F18").Value = 108
$wsExhibit.Range("F20").Value = 46
$wsExhibit.Range("F21").Value = 1579

# --- Sheet "演出" ---
$wsShow = $wb.Worksheets.Item("演出")
$wsShow.Range("F2").Value = 24

# --- Sheet "全部类型" ---
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F2").Value  = 24
$wsAll.Range("F3").Value  = 289
$wsAll.Range("F5").Value  = 10197
$wsAll.Range("F6").Value  = 330
$wsAll.Range("F7").Value  = 932
$wsAll.Range("F8").Value  = 1266
$wsAll.Range("F9").Value  = 6611
$wsAll.Range("F11").Value = 427
$wsAll.Range("F13").Value = 124
$wsAll.Range("F14").Value = 3145
$wsAll.Range("F17").Value = 619
$wsAll.Range("F18").Value = 119
$wsAll.Range("F19").Value = 108
$wsAll.Range("F21").Value = 46
$wsAll.Range("F22").Value = 1579
